# Apply the FA whitelist config edits:
#  - Add "T" markers to previously-empty cells (F2, F4, B8:B14)
#  - Remove the yellow highlight fill that was applied to G4 (style index 1)
#  - Move the active selection to G4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "T" markers added in column F for rows 2 and 4
$ws.Range("F2").Value = "T"
$ws.Range("F4").Value = "T"

# New "T" markers added in column B for rows 8-14
$ws.Range("B8").Value = "T"
$ws.Range("B9").Value = "T"
$ws.Range("B10").Value = "T"
$ws.Range("B11").Value = "T"
$ws.Range("B12").Value = "T"
$ws.Range("B13").Value = "T"
$ws.Range("B14").Value = "T"

# Remove the yellow fill highlight previously applied to G4
$ws.Range("G4").Interior.Pattern = -4142   # xlNone (no fill)

# Update the active selection to G4 (matches new cursor position in the diff)
$ws.Range("G4").Select()
